$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Elvis Galvis Galvis"
$ws.Range("E2").Value = 3.5
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = "Sí"

# Row 3
$ws.Range("B3").Value = "Albert Dayhan Diaz"
$ws.Range("F3").Value = 8
$ws.Range("I3").Value = 8

# Row 4
$ws.Range("B4").Value = "Duvan Gutierrez Lobo"

# Row 5
$ws.Range("B5").Value = "Daniela Guzman Perez"
$ws.Range("F5").Value = 7
$ws.Range("H5").Value = 8
$ws.Range("I5").Value = 15
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = "Sí"

# Row 6
$ws.Range("B6").Value = "Perez Carmen Andrea"

# Row 7
$ws.Range("B7").Value = "Angulo Juan Camilo"

# Row 8
$ws.Range("B8").Value = "Cristian Olivar Isaza"

# Row 9
$ws.Range("B9").Value = "Thania Milena Perez"

# Row 10
$ws.Range("B10").Value = "Marlene Ballena Guzman"
$ws.Range("J10").Value = 1

# Row 11
$ws.Range("B11").Value = "Jose Sierra Guzman"

# Row 12
$ws.Range("B12").Value = "Diego Silva Benavides"

# Row 13
$ws.Range("B13").Value = "Camilo Andres Daza"
$ws.Range("J13").Value = 1

# Row 14
$ws.Range("B14").Value = "Jose Castellano Endry"

# Row 15
$ws.Range("B15").Value = "Albert Dayhan Diaz"

# Row 16
$ws.Range("B16").Value = "Duvan Gutierrez Lobo"

# Row 17
$ws.Range("B17").Value = "Daniela Guzman Perez"

# Row 18
$ws.Range("B18").Value = "Perez Carmen Andrea"

# Row 19
$ws.Range("B19").Value = "Angulo Juan Camilo"

# Row 20
$ws.Range("B20").Value = "Cristian Olivar Isaza"

# Row 21
$ws.Range("B21").Value = "Thania Milena Perez"

# Row 22
$ws.Range("B22").Value = "Marlene Ballena Guzman"

# Row 23
$ws.Range("B23").Value = "Jose Sierra Guzman"

# Row 24
$ws.Range("B24").Value = "Diego Silva Benavides"
$ws.Range("J24").Value = 1

# Row 25
$ws.Range("B25").Value = "Camilo Andres Daza"

# Row 26
$ws.Range("B26").Value = "Jose Castellano Endry"

# Row 27
$ws.Range("B27").Value = "Albert Dayhan Diaz"

# Row 28
$ws.Range("B28").Value = "Duvan Gutierrez Lobo"
$ws.Range("J28").Value = 1

# Row 29
$ws.Range("B29").Value = "Daniela Guzman Perez"

# Row 30
$ws.Range("B30").Value = "Perez Carmen Andrea"
